$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2: existing "Test1" record becomes "Test2" (first/middle/last name
# and email change; address/zip/section/state/role/city/phones unchanged)
# ---------------------------------------------------------------------
$ws.Range("D2").Hyperlinks.Delete()

$ws.Range("A2").Value2 = "Test2"
$ws.Range("B2").Value2 = "TestMiddle2"
$ws.Range("C2").Value2 = "TestLast2"
$ws.Range("D2").Value2 = "test2@gmail.com"
$ws.Range("E2").Value2 = "Male"
$ws.Range("F2").Value2 = "11855 Lake Lucaya Drive Riverview Fl 33579"
$ws.Range("G2").Value2 = 33579
$ws.Range("H2").Value2 = 1
$ws.Range("I2").Value2 = "FL"
$ws.Range("J2").Value2 = 11
$ws.Range("K2").Value2 = "Riverview"
$ws.Range("L2").Value2 = "555-555-5555"
$ws.Range("M2").Value2 = "555-555-5555"
$ws.Range("N2").Value2 = "555-555-5555"

$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:test2@gmail.com", "", "", "test2@gmail.com") | Out-Null

# ---------------------------------------------------------------------
# Row 3: new record "Test3"
# ---------------------------------------------------------------------
$ws.Range("A3").Value2 = "Test3"
$ws.Range("B3").Value2 = "TestMiddle3"
$ws.Range("C3").Value2 = "TestLast3"
$ws.Range("D3").Value2 = "test3@gmail.com"
$ws.Range("E3").Value2 = "Male"
$ws.Range("F3").Value2 = "11856 Lake Lucaya Drive Riverview Fl 33579"
$ws.Range("G3").Value2 = 33579
$ws.Range("H3").Value2 = 0
$ws.Range("I3").Value2 = "FL"
$ws.Range("J3").Value2 = 11
$ws.Range("K3").Value2 = "Riverview"
$ws.Range("L3").Value2 = "555-555-5556"
$ws.Range("M3").Value2 = "555-555-5556"
$ws.Range("N3").Value2 = "555-555-5556"

$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:test3@gmail.com", "", "", "test3@gmail.com") | Out-Null

# ---------------------------------------------------------------------
# Row 4: new record "Test4"
# ---------------------------------------------------------------------
$ws.Range("A4").Value2 = "Test4"
$ws.Range("B4").Value2 = "TestMiddle4"
$ws.Range("C4").Value2 = "TestLast4"
$ws.Range("D4").Value2 = "test4@gmail.com"
$ws.Range("E4").Value2 = "Male"
$ws.Range("F4").Value2 = "11857 Lake Lucaya Drive Riverview Fl 33579"
$ws.Range("G4").Value2 = 33579
$ws.Range("H4").Value2 = 1
$ws.Range("I4").Value2 = "FL"
$ws.Range("J4").Value2 = 11
$ws.Range("K4").Value2 = "Riverview"
$ws.Range("L4").Value2 = "555-555-5557"
$ws.Range("M4").Value2 = "555-555-5557"
$ws.Range("N4").Value2 = "555-555-5557"

$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:test4@gmail.com", "", "", "test4@gmail.com") | Out-Null

# ---------------------------------------------------------------------
# Selection moved to K11 in the saved file
# ---------------------------------------------------------------------
$ws.Range("K11").Select() | Out-Null

Write-Host "done"
